$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# fix merge matrix bug: column B values for existing rows 1-44 were
# overstated (row 1 by 4, all others by 30); row 45 actually belongs to a
# later date with a different value, and two further merge-matrix data
# points (rows 46-47) were missing entirely.

$targets = @{
    1  = 250.0
    2  = 2737.0
    3  = 1538.0
    4  = 1894.0
    5  = 1560.0
    6  = 2205.0
    7  = 1940.0
    8  = 3542.0
    9  = 2775.0
    10 = 3553.0
    11 = 2342.0
    12 = 4492.0
    13 = 3997.0
    14 = 6608.0
    15 = 5547.0
    16 = 6547.0
    17 = 5345.0
    18 = 7445.0
    19 = 3575.0
    20 = 4924.0
    21 = 3812.0
    22 = 4505.0
    23 = 3215.0
    24 = 4375.0
    25 = 2730.0
    26 = 6080.0
    27 = 5010.0
    28 = 5511.0
    29 = 4578.0
    30 = 5192.0
    31 = 4387.0
    32 = 7765.0
    33 = 6858.0
    34 = 8558.0
    35 = 7902.0
    36 = 9355.0
    37 = 7394.0
    38 = 9120.0
    39 = 7930.0
    40 = 12252.0
    41 = 9740.0
    42 = 15227.0
    43 = 12990.0
    44 = 24939.0
}

foreach ($r in $targets.Keys) {
    $ws.Cells.Item($r, 2).Value = $targets[$r]
}

# Row 45 previously held 20200814 / 19230.0 — correct date + value.
# Column A holds the date as text, so force a text number format before
# assigning, otherwise the digit-only string gets stored as a number.
$ws.Cells.Item(45, 1).NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = "20201130"
$ws.Cells.Item(45, 2).Value = 12960.0

# Newly observed rows appended after row 45.
$ws.Cells.Item(46, 1).NumberFormat = "@"
$ws.Cells.Item(46, 1).Value = "20201228"
$ws.Cells.Item(46, 2).Value = 18520.0

$ws.Cells.Item(47, 1).NumberFormat = "@"
$ws.Cells.Item(47, 1).Value = "20210115"
$ws.Cells.Item(47, 2).Value = 14404.0
